# Put together BoM for the masthead light
# - Clear the AutoFilter criteria on column G (Part Number) so all rows show
#   (this also unhides the rows that the filter had hidden).
# - Move the active selection to F2:I2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the autofilter's filter criteria and unhide any rows the filter had
# hidden, while leaving the AutoFilter (the header-row dropdowns / the
# <autoFilter ref="A1:I55"/> range) itself in place.
$ws.ShowAllData()

# Move the selection to F2:I2 (active cell F2).
$ws.Range("F2:I2").Select()
